$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain text/number-safe cell updates (non-numeric-looking strings)
$ws.Range("D2").Value = '59.005.87'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '2.509.53'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  -0.12%  '
$ws.Range("E6").Value = '  -0.51%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  +0.46%  '
$ws.Range("E9").Value = '  +0.54%  '
$ws.Range("E10").Value = '  -1.08%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '2.954.64'
$ws.Range("E13").Value = '  +0.43%  '
$ws.Range("D14").Value = '58.953.21'
$ws.Range("E14").Value = '  +0.27%  '
$ws.Range("E15").Value = '  -1.51%  '
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '2.506.48'
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("E19").Value = '  +0.09%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("E21").Value = '  -0.04%  '
$ws.Range("E22").Value = '  +1.41%  '
$ws.Range("E23").Value = '  +0.58%  '
$ws.Range("E24").Value = '  +0.26%  '
$ws.Range("E25").Value = '  -0.56%  '
$ws.Range("E26").Value = '  +0.44%  '
$ws.Range("E27").Value = '  -0.70%  '
$ws.Range("D28").Value = '0.0₃0765'
$ws.Range("E28").Value = '  -0.66%  '
$ws.Range("E29").Value = '  -4.34%  '
$ws.Range("E30").Value = '  -0.94%  '
$ws.Range("E31").Value = '  +0.48%  '
$ws.Range("E32").Value = '  -0.01%  '
$ws.Range("E33").Value = '  -3.03%  '
$ws.Range("E34").Value = '  -0.75%  '
$ws.Range("E35").Value = '  -2.41%  '
$ws.Range("E36").Value = '  -1.11%  '
$ws.Range("E37").Value = '  -2.96%  '
$ws.Range("E38").Value = '  -1.09%  '
$ws.Range("E39").Value = '  -2.98%  '
$ws.Range("E40").Value = '  +0.82%  '
$ws.Range("E41").Value = '  -0.09%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("E42").Value = '  -3.90%  '
$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("E43").Value = '  +0.19%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("E45").Value = '  +1.49%  '
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("E47").Value = '  -2.74%  '
$ws.Range("E48").Value = '  -1.39%  '
$ws.Range("E49").Value = '  -0.12%  '
$ws.Range("D50").Value = '1.752.57'
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("E51").Value = '  -0.40%  '

# Numeric-looking text values in column D must be protected from Excel's
# automatic number coercion: format as Text, assign, then restore the
# default "Normal" style so no stray style index is left behind.
$protectedCells = @{
    "D5" = '533.71'
    "D6" = '136.03'
    "D11" = '5.41'
    "D12" = '0.346'
    "D15" = '22.78'
    "D20" = '324.11'
    "D22" = '5.94'
    "D23" = '65.06'
    "D24" = '0.420'
    "D26" = '1.00'
    "D27" = '7.53'
    "D29" = '6.45'
    "D31" = '168.64'
    "D34" = '18.39'
    "D36" = '4.05'
    "D38" = '3.58'
    "D39" = '0.800'
    "D40" = '281.69'
    "D42" = '5.01'
    "D43" = '0.604'
    "D44" = '10.94'
    "D45" = '129.82'
    "D49" = '17.30'
    "D51" = '0.983'
}
foreach ($ref in $protectedCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $protectedCells[$ref]
    $cell.Style = "Normal"
}
